# Append three new slides to the end of the deck -- "ID Selectors",
# "Class Selectors" and "Descendant Selectors" -- each built on the
# "Title and Content" layout (layout index 2) already used throughout
# this CSS-selectors deck (e.g. the preceding "Type Selectors" slide).
$p = $ppt.ActivePresentation

# ---- New slide 29: "ID Selectors" ----
$slide = $p.Slides.Add(29, 2)
$tr = $slide.Shapes.Item(1).TextFrame.TextRange  # title placeholder
$tr.Text = "ID "
[void]$tr.InsertAfter("Selectors")
$tr = $slide.Shapes.Item(2).TextFrame.TextRange  # body placeholder
$tr.Text = "A"
[void]$tr.InsertAfter("n ")
[void]$tr.InsertAfter("ID selector, lets us assign a unique ID to an element.")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("That ")
[void]$tr.InsertAfter("way, we're able to specifically target an element based ")
[void]$tr.InsertAfter("on its ")
[void]$tr.InsertAfter("ID attribute")
[void]$tr.InsertAfter(".")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("ID selectors are declared using the pound ")
[void]$tr.InsertAfter("(#) symbol ")
[void]$tr.InsertAfter("followed by the ID name")
[void]$tr.InsertAfter(".")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("ID's are unique to the page, ")
[void]$tr.InsertAfter("so it's ")
[void]$tr.InsertAfter("important to remember that an element can only have one ID.")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("A ")
[void]$tr.InsertAfter("page can only have one element with the same ID name")
[void]$tr.InsertAfter(".")
$slide.Shapes.Item(2).TextFrame.AutoSize = 2  # ppAutoSizeTextToFitShape -> normAutofit

# ---- New slide 30: "Class Selectors" ----
$slide = $p.Slides.Add(30, 2)
$tr = $slide.Shapes.Item(1).TextFrame.TextRange  # title placeholder
$tr.Text = "Class "
[void]$tr.InsertAfter("Selectors")
$tr = $slide.Shapes.Item(2).TextFrame.TextRange  # body placeholder
$tr.Text = "Class selectors let us target elements based on their class attribute. The main difference between a class and an ID selector is that IDs are unique and they’re used to identify one element on the page, whereas a class can target more than one element"
[void]$tr.InsertAfter(".")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("T")
[void]$tr.InsertAfter("he ")
[void]$tr.InsertAfter("main difference between a class ")
[void]$tr.InsertAfter("and an ")
[void]$tr.InsertAfter("ID selector is that ID's are ")
[void]$tr.InsertAfter("unique. And ")
[void]$tr.InsertAfter("they're used to identify one element on the page.")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("Whereas ")
[void]$tr.InsertAfter("a class can be used to classify and target more than one ")
[void]$tr.InsertAfter("element. This ")
[void]$tr.InsertAfter("makes classes more flexible than ")
[void]$tr.InsertAfter("ID's.")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("Classes let us target more than one element with the same class ")
[void]$tr.InsertAfter("name. In ")
[void]$tr.InsertAfter("fact, that's one of the biggest advantages to using class selectors.")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("Multiple ")
[void]$tr.InsertAfter("elements can share the same class, ")
[void]$tr.InsertAfter("so we're ")
[void]$tr.InsertAfter("able to reuse them throughout a page")
[void]$tr.InsertAfter(".")
$slide.Shapes.Item(2).TextFrame.AutoSize = 2  # ppAutoSizeTextToFitShape -> normAutofit

# ---- New slide 31: "Descendant Selectors" ----
$slide = $p.Slides.Add(31, 2)
$tr = $slide.Shapes.Item(1).TextFrame.TextRange  # title placeholder
$tr.Text = "Descendant "
[void]$tr.InsertAfter("Selectors")
$tr = $slide.Shapes.Item(2).TextFrame.TextRange  # body placeholder
$tr.Text = "CSS also lets us target elements based on their relationship in the HTML document."
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("For ")
[void]$tr.InsertAfter("instance, we can combine selectors to create what's called a ")
[void]$tr.InsertAfter("descendent selector")
[void]$tr.InsertAfter(", because it targets an element that's a descendent of another element.")
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter("This ")
[void]$tr.InsertAfter("makes our selectors more specific")
[void]$tr.InsertAfter(".")
